$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1867.3334
$ws.Range("J32").Value = 2002
$ws.Range("L32").Value = 2002
$ws.Range("N32").Value = -2654
$ws.Range("H40").Value = 1826.9131
$ws.Range("I40").Value = 1741.2354
$ws.Range("J40").Value = 2069.6667
$ws.Range("K40").Value = 1741.2354
$ws.Range("L40").Value = 2069.6667
$ws.Range("M40").Value = -1566.2354
$ws.Range("N40").Value = -2419.6667
$ws.Range("H74").Value = 5964.3335
$ws.Range("I74").Value = 4796
$ws.Range("K74").Value = 4796
$ws.Range("M74").Value = -3860
$ws.Range("H77").Value = 5964.3335
$ws.Range("I77").Value = 4796
$ws.Range("K77").Value = 23980
$ws.Range("M77").Value = -19300
$ws.Range("H96").Value = 778.4
$ws.Range("I96").Value = 651.5
$ws.Range("J96").Value = 968.75
$ws.Range("K96").Value = 1954.5
$ws.Range("L96").Value = 2906.25
$ws.Range("M96").Value = -581.5
$ws.Range("N96").Value = -5652.25
$ws.Range("H98").Value = 876.61536
$ws.Range("I98").Value = 860.96
$ws.Range("J98").Value = 904.5714
$ws.Range("K98").Value = 860.96
$ws.Range("L98").Value = 904.5714
$ws.Range("M98").Value = 637.04
$ws.Range("N98").Value = -3900.5714
$ws.Range("H100").Value = 1007.3077
$ws.Range("I100").Value = 1227.1428
$ws.Range("K100").Value = 1227.1428
$ws.Range("M100").Value = -686.1428000000001
$ws.Range("H106").Value = 1611.75
$ws.Range("I106").Value = 1659.8667
$ws.Range("K106").Value = 1659.8667
$ws.Range("M106").Value = -1028.8667
$ws.Range("H116").Value = 1788.8667
$ws.Range("I116").Value = 1144.1538
$ws.Range("K116").Value = 1144.1538
$ws.Range("M116").Value = 2297.8462
$ws.Range("H122").Value = 876.61536
$ws.Range("I122").Value = 860.96
$ws.Range("J122").Value = 904.5714
$ws.Range("K122").Value = 2582.88
$ws.Range("L122").Value = 2713.7142
$ws.Range("M122").Value = -132.8800000000001
$ws.Range("N122").Value = -7613.7142
$ws.Range("H129").Value = 917.13336
$ws.Range("J129").Value = 971.80554
$ws.Range("L129").Value = 2915.41662
$ws.Range("N129").Value = -12915.41662
$ws.Range("H141").Value = 2353.75
$ws.Range("I141").Value = 1309.091
$ws.Range("J141").Value = 6184.1665
$ws.Range("K141").Value = 3927.273
$ws.Range("L141").Value = 18552.4995
$ws.Range("M141").Value = 1252.727
$ws.Range("N141").Value = -28912.4995

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 60351
$ws.Range("I2").Value = 1165.6666
$ws.Range("J2").Value = 202395.8
$ws.Range("K2").Value = 1165.6666
$ws.Range("L2").Value = 202395.8
$ws.Range("M2").Value = -1052.6666
$ws.Range("N2").Value = -202621.8
$ws.Range("H32").Value = 32234.197
$ws.Range("I32").Value = 5964.35
$ws.Range("K32").Value = 5964.35
$ws.Range("M32").Value = -5677.35
$ws.Range("H92").Value = 13550
$ws.Range("J92").Value = 13550
$ws.Range("L92").Value = 13550
$ws.Range("N92").Value = -18542
$ws.Range("H97").Value = 45266.566
$ws.Range("I97").Value = 67179.87
$ws.Range("J97").Value = 4179.125
$ws.Range("K97").Value = 67179.87
$ws.Range("L97").Value = 4179.125
$ws.Range("M97").Value = -66683.87
$ws.Range("N97").Value = -5171.125
$ws.Range("H116").Value = 60351
$ws.Range("I116").Value = 1165.6666
$ws.Range("J116").Value = 202395.8
$ws.Range("K116").Value = 1165.6666
$ws.Range("L116").Value = 202395.8
$ws.Range("M116").Value = 1128.3334
$ws.Range("N116").Value = -206983.8
$ws.Range("H119").Value = 44000
$ws.Range("J119").Value = 44000
$ws.Range("L119").Value = 44000
$ws.Range("N119").Value = -53676
$ws.Range("H122").Value = 1968.6129
$ws.Range("I122").Value = 1781.375
$ws.Range("J122").Value = 2610.5715
$ws.Range("K122").Value = 5344.125
$ws.Range("L122").Value = 7831.7145
$ws.Range("M122").Value = -2894.125
$ws.Range("N122").Value = -12731.7145

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 60351
$ws.Range("I3").Value = 1165.6666
$ws.Range("J3").Value = 202395.8
$ws.Range("K3").Value = 1165.6666
$ws.Range("L3").Value = 202395.8
$ws.Range("M3").Value = -1051.6666
$ws.Range("N3").Value = -202623.8
$ws.Range("H22").Value = 376.33334
$ws.Range("I22").Value = 371.6
$ws.Range("K22").Value = 371.6
$ws.Range("M22").Value = -198.6
$ws.Range("H35").Value = 19218
$ws.Range("I35").Value = 10000
$ws.Range("J35").Value = 20370.25
$ws.Range("K35").Value = 10000
$ws.Range("L35").Value = 20370.25
$ws.Range("M35").Value = -9690
$ws.Range("N35").Value = -20990.25
$ws.Range("H99").Value = 1538.6875
$ws.Range("I99").Value = 1190.8889
$ws.Range("J99").Value = 1985.8572
$ws.Range("K99").Value = 1190.8889
$ws.Range("L99").Value = 1985.8572
$ws.Range("M99").Value = 307.1111000000001
$ws.Range("N99").Value = -4981.8572
$ws.Range("H105").Value = 335948.5
$ws.Range("I105").Value = 335593.34
$ws.Range("J105").Value = 336303.66
$ws.Range("K105").Value = 335593.34
$ws.Range("L105").Value = 336303.66
$ws.Range("M105").Value = -333846.34
$ws.Range("N105").Value = -339797.66
$ws.Range("H134").Value = 2111
$ws.Range("I134").Value = 2234.35
$ws.Range("K134").Value = 6703.049999999999
$ws.Range("M134").Value = -4168.049999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15741.357
$ws.Range("I31").Value = 26597.05
$ws.Range("J31").Value = 2084.1936
$ws.Range("K31").Value = 26597.05
$ws.Range("L31").Value = 2084.1936
$ws.Range("M31").Value = -26302.05
$ws.Range("N31").Value = -2674.1936
$ws.Range("H34").Value = 15741.357
$ws.Range("I34").Value = 26597.05
$ws.Range("J34").Value = 2084.1936
$ws.Range("K34").Value = 26597.05
$ws.Range("L34").Value = 2084.1936
$ws.Range("M34").Value = -26395.05
$ws.Range("N34").Value = -2488.1936
$ws.Range("H41").Value = 10261.583
$ws.Range("J41").Value = 11189.091
$ws.Range("L41").Value = 11189.091
$ws.Range("N41").Value = -12045.091
$ws.Range("H50").Value = 10594.286
$ws.Range("J50").Value = 10594.286
$ws.Range("L50").Value = 10594.286
$ws.Range("N50").Value = -11844.286
$ws.Range("H51").Value = 6992.25
$ws.Range("J51").Value = 7966.3335
$ws.Range("L51").Value = 7966.3335
$ws.Range("N51").Value = -9438.333500000001
$ws.Range("H59").Value = 24805
$ws.Range("J59").Value = 24805
$ws.Range("L59").Value = 24805
$ws.Range("N59").Value = -27095
$ws.Range("H61").Value = 6992.25
$ws.Range("J61").Value = 7966.3335
$ws.Range("L61").Value = 7966.3335
$ws.Range("N61").Value = -8662.333500000001
$ws.Range("H86").Value = 2429.963
$ws.Range("I86").Value = 1899.2307
$ws.Range("K86").Value = 1899.2307
$ws.Range("M86").Value = -776.2307000000001
$ws.Range("H89").Value = 2429.963
$ws.Range("I89").Value = 1899.2307
$ws.Range("K89").Value = 9496.1535
$ws.Range("M89").Value = -3880.1535

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1992
$ws.Range("I132").Value = 1656.1428
$ws.Range("K132").Value = 14905.2852
$ws.Range("M132").Value = -12375.2852

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 5370.5386
$ws.Range("I43").Value = 1933.3334
$ws.Range("J43").Value = 6401.7
$ws.Range("K43").Value = 1933.3334
$ws.Range("L43").Value = 6401.7
$ws.Range("M43").Value = -1782.3334
$ws.Range("N43").Value = -6703.7
$ws.Range("H132").Value = 3200.55
$ws.Range("I132").Value = 2900.8
$ws.Range("J132").Value = 4099.8
$ws.Range("K132").Value = 8702.400000000001
$ws.Range("L132").Value = 12299.4
$ws.Range("M132").Value = -6172.400000000001
$ws.Range("N132").Value = -17359.4
$ws.Range("H135").Value = 39974.316
$ws.Range("J135").Value = 39974.316
$ws.Range("L135").Value = 39974.316
$ws.Range("N135").Value = -50114.316

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2461.1538
$ws.Range("I46").Value = 699
$ws.Range("J46").Value = 3562.5
$ws.Range("K46").Value = 699
$ws.Range("L46").Value = 3562.5
$ws.Range("M46").Value = -511
$ws.Range("N46").Value = -3938.5
$ws.Range("H93").Value = 1641.8889
$ws.Range("I93").Value = 1659.625
$ws.Range("J93").Value = 1500
$ws.Range("K93").Value = 1659.625
$ws.Range("L93").Value = 1500
$ws.Range("M93").Value = -411.625
$ws.Range("N93").Value = -3996
$ws.Range("H94").Value = 30000
$ws.Range("J94").Value = 30000
$ws.Range("L94").Value = 30000
$ws.Range("N94").Value = -31352
$ws.Range("H132").Value = 7192.45
$ws.Range("I132").Value = 7911.615
$ws.Range("J132").Value = 5856.857
$ws.Range("K132").Value = 23734.845
$ws.Range("L132").Value = 17570.571
$ws.Range("M132").Value = -21204.845
$ws.Range("N132").Value = -22630.571

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3765.2917
$ws.Range("I132").Value = 5221.6665
$ws.Range("J132").Value = 2308.9167
$ws.Range("K132").Value = 15664.9995
$ws.Range("L132").Value = 6926.750100000001
$ws.Range("M132").Value = -13134.9995
$ws.Range("N132").Value = -11986.7501
$ws.Range("H136").Value = 726.7436
$ws.Range("I136").Value = 572.6818
$ws.Range("J136").Value = 926.1177
$ws.Range("K136").Value = 1718.0454
$ws.Range("L136").Value = 2778.3531
$ws.Range("M136").Value = 831.9546
$ws.Range("N136").Value = -7878.3531
